$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

$ws.Range("B2").Value = 14.718302893770216
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 15.195459151301188
$ws.Range("E2").ClearContents()

$ws.Range("B3").Value = 10.457156071013809
$ws.Range("C3").Value = -12.183051192106124
$ws.Range("D3").Value = 10.067551391207463
$ws.Range("E3").Value = -16.548441345687586

$ws.Range("B1:E3").Select()
